# refresh data july 5th
# Insert a new first column ("Name") in the TERI table, re-populate the
# table (rows re-ordered + "Training ?" status line content added), fix
# mis-encoded punctuation in the two long description cells, widen/resize
# a few columns and tweak two alignment styles.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a brand-new column A - everything else (B:E) shifts right to
#    (C:F) and all formatting / merged cells shift with it automatically.
$ws.Columns.Item(1).Insert()

# 2) Header row (row 2) - new "Name" header in col A, existing headers now
#    live in B2:F2 (already shifted correctly by the Insert above).
$ws.Range("A2").Value2 = "Name"

# 3) Row 3 ("Outputs" banner) - A3 already holds "Outputs" after the shift;
#    nothing else to do there, the merge already grew to A3:F3.

# 4) Data rows 4-6. The underlying records were also re-ordered, so we
#    rewrite B4:F6 (old A:E data) in full alongside the new A4:A6 "Name"
#    values rather than trying to reconcile a partial shuffle.

# Row 4 -> "5 articles" / Communication pieces record
$ws.Range("A4").Value2 = "5 articles"
$ws.Range("B4").Value2 = "In Process"
$ws.Range("C4").Value2 = "2022-First Half"
$ws.Range("D4").Value2 = "Communication pieces developed"
$ws.Range("E4").Value2 = "Improve public awareness and usage of public transport"
$ws.Range("F4").Value2 = "5 articles - the research team has already published articles in news media and magazines: 1. Money Control: Net-Zero Emissions | We need diverse strategies to decarbonize India’s`n"

# Row 5 -> "2 policy briefs on private freight" / Policy record
$ws.Range("A5").Value2 = "2 policy briefs on private freight"
$ws.Range("B5").Value2 = "Completed"
$ws.Range("C5").Value2 = "2022-First Half"
$ws.Range("D5").Value2 = "Policy and regulatory recommendations"
$ws.Range("E5").Value2 = "Increased share of railways in freight transport"
$ws.Range("F5").Value2 = "2 policy briefs published, focused on private freight terminals and freight forwarder scheme of IR namely – ‘Developing Rail Freight Terminals: Energizing Private Partnerships’ and  Moving Towards Aggregation: Freight Forwarders Scheme of Railways’.https://www.teriin.org/policy-brief/developing-rail-freight-terminals-energizing-private-partnerships https://www.teriin.org/policy-brief/moving-towards-aggregation-freight-forwarder-scheme-railways"

# Row 6 -> "Training ?" / Trainings record
$ws.Range("A6").Value2 = "Training ?"
$ws.Range("B6").Value2 = "In Process"
$ws.Range("C6").Value2 = "2022-First Half"
$ws.Range("D6").Value2 = "Trainings/Webinars/Seminars"
$ws.Range("E6").Value2 = "Increased share of railways in freight transport"
$ws.Range("F6").Value2 = "Not done yet"

# 5) Column widths.
$ws.Columns.Item(1).ColumnWidth = 15.71
$ws.Columns.Item(3).ColumnWidth = 13.71
$ws.Columns.Item(5).ColumnWidth = 30.71
$ws.Columns.Item(6).ColumnWidth = 70.71

# 6) Style tweaks: title cell (A1) now centers vertically, and the
#    "Outputs" banner row (A3:F3) gains horizontal centering.
$ws.Range("A1").VerticalAlignment = -4108
$ws.Range("A3:F3").HorizontalAlignment = -4108
